$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the referenced table entry from Buffsht/Buff/item.xlsx to
# UnitTable/UnitTableData/unit.xlsx
$ws.Range("B4").Value = "UnitTable"
$ws.Range("C4").Value = "UnitTableData"

# Re-point the "input" cell (E4) to the new unit.xlsx workbook and turn it
# into a hyperlink, the way item.xlsx used to be referenced.
$ws.Hyperlinks.Add($ws.Range("E4"), "unit.xlsx", "", "", "UnitTable@unit.xlsx")

# Update the active selection left behind by the editing session.
$ws.Range("E7").Select() | Out-Null
